# Add a new publication row to the "Journals" worksheet.
# A new row is inserted at row 3 (pushing all existing entries down by one)
# for the new "In Press" publication:
#   Barney, J. L., Barrett, T. S., Lensegrav-Benson, T., Quakenbush, B., Twohig, M. P.
#   Examining a mediation model of body image-related cognitive fusion, intuitive
#   eating, and eating disorder symptom severity in a clinical sample.
#   Eating and Weight Disorders - Studies on Anorexia, Bulimia, and Obesity

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journals")

# Shift rows 3:53 down to make room for the new entry.
$ws.Rows.Item(3).EntireRow.Insert()

$ws.Range("A3").Value = "Barney, J. L., Barrett, T. S., Lensegrav-Benson, T., Quakenbush, B., Twohig, M. P."
$ws.Range("B3").Value = "In Press"
$ws.Range("C3").Value = "Examining a mediation model of body image-related cognitive fusion, intuitive eating, and eating disorder symptom severity in a clinical sample."
$ws.Range("D3").Value = [char]0x2013
$ws.Range("D3").Value = "Eating and Weight Disorders " + [char]0x2013 + " Studies on Anorexia, Bulimia, and Obesity"

$ws.Rows.Item(3).RowHeight = 51

$ws.Range("C5").Select()
